# Apply the "end of loop" progress update: refresh quantity / amount columns
# on the Bill Summary sheet with the latest computed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple numeric "Qty executed upto date" (column C) updates ---------
$ws.Range("C8").Value  = 67
$ws.Range("C9").Value  = 65
$ws.Range("C10").Value = 54
$ws.Range("C11").Value = 14
$ws.Range("C12").Value = 58
$ws.Range("C13").Value = 80
$ws.Range("C14").Value = 67
$ws.Range("C15").Value = 54
$ws.Range("C16").Value = 50
$ws.Range("C17").Value = 83

# --- "Upto date Amount" / total cells (column G & H) -------------------
# These cells hold text-formatted numbers (e.g. "16640.00") rather than
# real numeric values, so assigning a plain string would make Excel
# auto-convert it to a number and drop the formatting. Instead we write
# a literal-string formula and immediately flatten it back down to a
# plain cached value with Copy/PasteSpecial(values), which keeps the
# text representation intact without touching any cell styles.
$textAmounts = [ordered]@{
    "G9"  = "16640.00"
    "G10" = "25488.00"
    "G11" = "9268.00"
    "G13" = "10880.00"
    "G14" = "1541.00"
    "G19" = "63817.00"
    "H19" = "63817.00"
    "G21" = "63817.00"
    "H21" = "63817.00"
}

foreach ($addr in $textAmounts.Keys) {
    $ws.Range($addr).Formula = "=""" + $textAmounts[$addr] + """"
}
foreach ($addr in $textAmounts.Keys) {
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
